# Generate Report for Handback
# Update timestamp cells to reflect newly generated handback report times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for 65ba8943... row (shared between Overview!G2 and de-de!H2)
$wsOverview.Range("G2").Value = "2016-08-16 13:05:39"
$wsDeDe.Range("H2").Value = "2016-08-16 13:05:39"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 65ba8943... row
$wsZhCn.Range("H2").Value = "2016-08-16 13:05:33"
$wsZhCn.Range("K2").Value = "2016-08-16 13:05:50"

# de-de sheet: Correspond Handback DateTime for 65ba8943... row
$wsDeDe.Range("K2").Value = "2016-08-16 13:05:58"
